# Tasks.docx edit:
#  - Remove the struck-through "Über achtzig Vorstandsaufgaben (Alternativer
#    Titel)" paragraph entirely (its paragraph mark too, so the following
#    "Einblick in die Vorstandsarbeit" paragraph slides up and keeps its own,
#    non-struck-through paragraph formatting).
#  - Move the (hidden) "_GoBack" bookmark from its old spot, in the middle of
#    the last real paragraph ("... gibt immer wieder |neue Versionen."), to
#    the very start of the paragraph that used to hold the struck-through
#    title (now the start of the "Einblick ..." paragraph). Re-adding a
#    "_GoBack" bookmark anywhere automatically relocates the single allowed
#    instance, so this also takes care of deleting it from its old location.

$d = $word.ActiveDocument

# Locate the struck-through title paragraph by its text and delete it
# (delete the whole paragraph range, including the paragraph mark, so the
# next paragraph absorbs its position without picking up the strike
# formatting).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Über achtzig Vorstandsaufgaben*Alternativer Titel*") {
        $target = $p
        break
    }
}
$target.Range.Delete()

# Re-find the paragraph that now starts with "Einblick" (previously the
# second of the pair, now holding the merged content) and drop the
# "_GoBack" bookmark right at its start.
$einblick = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Einblick*") {
        $einblick = $p
        break
    }
}
$bmStart = $d.Range($einblick.Range.Start, $einblick.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmStart)
